$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.6084629976340921
$ws.Cells.Item(2, 3).Value = 0.06128436316929964
$ws.Cells.Item(2, 5).Value = 0.08626289677362209
$ws.Cells.Item(2, 6).Value = 0.4443680307746263
$ws.Cells.Item(2, 7).Value = 1.185236305379831
$ws.Cells.Item(2, 8).Value = 1.163481456830226
$ws.Cells.Item(2, 9).Value = 1.154009388848159
$ws.Cells.Item(2, 11).Value = 0.390417530626479
$ws.Cells.Item(2, 12).Value = 0.2173515428299595

$ws.Cells.Item(3, 2).Value = 0.5693292127438667
$ws.Cells.Item(3, 3).Value = 0.0590161716693558
$ws.Cells.Item(3, 5).Value = 0.08563889364818777
$ws.Cells.Item(3, 6).Value = 0.387822817061874
$ws.Cells.Item(3, 7).Value = 1.192172659259214
$ws.Cells.Item(3, 8).Value = 1.171935428554775
$ws.Cells.Item(3, 9).Value = 1.163385801421644
$ws.Cells.Item(3, 11).Value = 0.35419849425449
$ws.Cells.Item(3, 12).Value = 0.2102882849231094

$ws.Cells.Item(4, 2).Value = 0.5455213825979683
$ws.Cells.Item(4, 3).Value = 0.05760419518581017
$ws.Cells.Item(4, 5).Value = 0.0852968503746574
$ws.Cells.Item(4, 6).Value = 0.3531389305168915
$ws.Cells.Item(4, 7).Value = 1.197170276363337
$ws.Cells.Item(4, 8).Value = 1.177644995680595
$ws.Cells.Item(4, 9).Value = 1.169727700752176
$ws.Cells.Item(4, 11).Value = 0.3320555347040965
$ws.Cells.Item(4, 12).Value = 0.2060716361910551

$ws.Cells.Item(5, 2).Value = 0.5358754469863811
$ws.Cells.Item(5, 3).Value = 0.05702394316969617
$ws.Cells.Item(5, 5).Value = 0.08516781551899832
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 1.199392315398313
$ws.Cells.Item(5, 8).Value = 1.180102126618863
$ws.Cells.Item(5, 9).Value = 1.172459016779648
$ws.Cells.Item(5, 11).Value = 0.3230565135253585
$ws.Cells.Item(5, 12).Value = 0.2043835887605923

$ws.Cells.Item(6, 2).Value = 0.5342771409589204
$ws.Cells.Item(6, 3).Value = 0.05692729900973603
$ws.Cells.Item(6, 5).Value = 0.08514701501025357
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 1.199772477150276
$ws.Cells.Item(6, 8).Value = 1.180518008807184
$ws.Cells.Item(6, 9).Value = 1.172921421706793
$ws.Cells.Item(6, 11).Value = 0.3215637180787354
$ws.Cells.Item(6, 12).Value = 0.2041051192945531

$ws.Cells.Item(7, 2).Value = 0.5453910669257311
$ws.Cells.Item(7, 3).Value = 0.0575963893927991
$ws.Cells.Item(7, 5).Value = 0.08529506823448685
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 1.197199492967059
$ws.Cells.Item(7, 8).Value = 1.17767760530576
$ws.Cells.Item(7, 9).Value = 1.169763941411585
$ws.Cells.Item(7, 11).Value = 0.3319340714280656
$ws.Cells.Item(7, 12).Value = 0.2060487479709536

$ws.Cells.Item(8, 2).Value = 0.5949242006181237
$ws.Cells.Item(8, 3).Value = 0.06050629080995407
$ws.Cells.Item(8, 5).Value = 0.08603922016225241
$ws.Cells.Item(8, 6).Value = 0.4248636149813336
$ws.Cells.Item(8, 7).Value = 1.1874744833367
$ws.Cells.Item(8, 8).Value = 1.166288698749952
$ws.Cells.Item(8, 9).Value = 1.157120953129926
$ws.Cells.Item(8, 11).Value = 0.3779095506611441
$ws.Cells.Item(8, 12).Value = 0.2148912014448854

$ws.Cells.Item(9, 2).Value = 0.6937910124604798
$ws.Cells.Item(9, 3).Value = 0.06606006013595334
$ws.Cells.Item(9, 5).Value = 0.08782404619279305
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 1.174277872801255
$ws.Cells.Item(9, 8).Value = 1.148072590491495
$ws.Cells.Item(9, 9).Value = 1.136972405237181
$ws.Cells.Item(9, 11).Value = 0.4688156514103241
$ws.Cells.Item(9, 12).Value = 0.2331845079563664

$ws.Cells.Item(10, 2).Value = 0.7674700932698784
$ws.Cells.Item(10, 3).Value = 0.07004857036520207
$ws.Cells.Item(10, 5).Value = 0.08933335343022364
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 1.168182523059812
$ws.Cells.Item(10, 8).Value = 1.137201239682355
$ws.Cells.Item(10, 9).Value = 1.125007280633397
$ws.Cells.Item(10, 11).Value = 0.5360531758874458
$ws.Cells.Item(10, 12).Value = 0.2472066757263605

$ws.Cells.Item(11, 2).Value = 0.8012122772076395
$ws.Cells.Item(11, 3).Value = 0.07184337261864471
$ws.Cells.Item(11, 5).Value = 0.09006288811521657
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 1.166195543792156
$ws.Cells.Item(11, 8).Value = 1.132801498720781
$ws.Cells.Item(11, 9).Value = 1.120181769133779
$ws.Cells.Item(11, 11).Value = 0.5667376476852155
$ws.Cells.Item(11, 12).Value = 0.2537124279309353

$ws.Cells.Item(12, 2).Value = 0.8140215634044523
$ws.Cells.Item(12, 3).Value = 0.07252021430925026
$ws.Cells.Item(12, 5).Value = 0.09034530703496202
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 1.16555645358541
$ws.Cells.Item(12, 8).Value = 1.131213940935979
$ws.Cells.Item(12, 9).Value = 1.118443390710375
$ws.Cells.Item(12, 11).Value = 0.5783708744558851
$ws.Cells.Item(12, 12).Value = 0.2561942368464969

$ws.Cells.Item(13, 2).Value = 0.811261446142197
$ws.Cells.Item(13, 3).Value = 0.07237456925913932
$ws.Cells.Item(13, 5).Value = 0.09028420934336268
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 1.165689047129788
$ws.Cells.Item(13, 8).Value = 1.1315523563825
$ws.Cells.Item(13, 9).Value = 1.118813823757414
$ws.Cells.Item(13, 11).Value = 0.5758648483445938
$ws.Cells.Item(13, 12).Value = 0.2556589251153554

$ws.Cells.Item(14, 2).Value = 0.8022654695242011
$ws.Cells.Item(14, 3).Value = 0.07189911309127694
$ws.Cells.Item(14, 5).Value = 0.09008599952183971
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 1.16614069212288
$ws.Cells.Item(14, 8).Value = 1.132669315216233
$ws.Cells.Item(14, 9).Value = 1.120036968414468
$ws.Cells.Item(14, 11).Value = 0.5676944481744215
$ws.Cells.Item(14, 12).Value = 0.2539162428833208

$ws.Cells.Item(15, 2).Value = 0.7967593066862548
$ws.Cells.Item(15, 3).Value = 0.07160751657507092
$ws.Cells.Item(15, 5).Value = 0.08996539204141385
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 1.166432107311152
$ws.Cells.Item(15, 8).Value = 1.133363713257708
$ws.Cells.Item(15, 9).Value = 1.12079776649076
$ws.Cells.Item(15, 11).Value = 0.5626916129522215
$ws.Cells.Item(15, 12).Value = 0.252851171512134

$ws.Cells.Item(16, 2).Value = 0.7652694483964808
$ws.Cells.Item(16, 3).Value = 0.06993088212145437
$ws.Cells.Item(16, 5).Value = 0.08928653944910891
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 1.168328218291848
$ws.Cells.Item(16, 8).Value = 1.137499758822969
$ws.Cells.Item(16, 9).Value = 1.125335076869362
$ws.Cells.Item(16, 11).Value = 0.5340498087755066
$ws.Cells.Item(16, 12).Value = 0.2467840615839947

$ws.Cells.Item(17, 2).Value = 0.746008741064486
$ws.Cells.Item(17, 3).Value = 0.06889730780996217
$ws.Cells.Item(17, 5).Value = 0.08888107332995077
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 1.169692928252047
$ws.Cells.Item(17, 8).Value = 1.140176895115303
$ws.Cells.Item(17, 9).Value = 1.128276825072831
$ws.Cells.Item(17, 11).Value = 0.5165037645626853
$ws.Cells.Item(17, 12).Value = 0.243094584547606

$ws.Cells.Item(18, 2).Value = 0.7349517200042044
$ws.Cells.Item(18, 3).Value = 0.068300978820659
$ws.Cells.Item(18, 5).Value = 0.08865190287335878
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 1.170551818972569
$ws.Cells.Item(18, 8).Value = 1.141768072554257
$ws.Cells.Item(18, 9).Value = 1.130026952650915
$ws.Cells.Item(18, 11).Value = 0.5064209703550944
$ws.Cells.Item(18, 12).Value = 0.240984450974679

$ws.Cells.Item(19, 2).Value = 0.7312116648722338
$ws.Cells.Item(19, 3).Value = 0.06809875511106611
$ws.Cells.Item(19, 5).Value = 0.08857500453586908
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 1.17085531464258
$ws.Cells.Item(19, 8).Value = 1.142315637149281
$ws.Cells.Item(19, 9).Value = 1.130629492138411
$ws.Cells.Item(19, 11).Value = 0.5030087056903199
$ws.Cells.Item(19, 12).Value = 0.240272050697186

$ws.Cells.Item(20, 2).Value = 0.7480568830374636
$ws.Cells.Item(20, 3).Value = 0.06900752440244418
$ws.Cells.Item(20, 5).Value = 0.08892381757287282
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 1.169539997019456
$ws.Cells.Item(20, 8).Value = 1.139886593447329
$ws.Cells.Item(20, 9).Value = 1.127957655549196
$ws.Cells.Item(20, 11).Value = 0.5183706187802954
$ws.Cells.Item(20, 12).Value = 0.2434860985654126

$ws.Cells.Item(21, 2).Value = 0.8049069456337179
$ws.Cells.Item(21, 3).Value = 0.07203884233253177
$ws.Cells.Item(21, 5).Value = 0.0901440514926044
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 1.166004954474147
$ws.Cells.Item(21, 8).Value = 1.132339105613866
$ws.Cells.Item(21, 9).Value = 1.11967528627202
$ws.Cells.Item(21, 11).Value = 0.5700939243953087
$ws.Cells.Item(21, 12).Value = 0.2544276166860726

$ws.Cells.Item(22, 2).Value = 0.8422471710227626
$ws.Cells.Item(22, 3).Value = 0.07400360234400694
$ws.Cells.Item(22, 5).Value = 0.09097744286129483
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 1.164355362292653
$ws.Cells.Item(22, 8).Value = 1.127864138756166
$ws.Cells.Item(22, 9).Value = 1.114780745492475
$ws.Cells.Item(22, 11).Value = 0.6039776779192323
$ws.Cells.Item(22, 12).Value = 0.2616847229710686

$ws.Cells.Item(23, 2).Value = 0.8223011978688817
$ws.Cells.Item(23, 3).Value = 0.07295646976369596
$ws.Cells.Item(23, 5).Value = 0.09052936655799115
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 1.165175210314032
$ws.Cells.Item(23, 8).Value = 1.130210608755007
$ws.Cells.Item(23, 9).Value = 1.117345566471592
$ws.Cells.Item(23, 11).Value = 0.5858861231597814
$ws.Cells.Item(23, 12).Value = 0.2578017655150973

$ws.Cells.Item(24, 2).Value = 0.7471308675012551
$ws.Cells.Item(24, 3).Value = 0.06895770206501339
$ws.Cells.Item(24, 5).Value = 0.08890448063530343
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 1.169608905808417
$ws.Cells.Item(24, 8).Value = 1.140017676679065
$ws.Cells.Item(24, 9).Value = 1.128101768697732
$ws.Cells.Item(24, 11).Value = 0.5175265994453468
$ws.Cells.Item(24, 12).Value = 0.2433090608256521

$ws.Cells.Item(25, 2).Value = 0.6668609646306436
$ws.Cells.Item(25, 3).Value = 0.06457383069820821
$ws.Cells.Item(25, 5).Value = 0.08730639983112454
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 1.177217057977501
$ws.Cells.Item(25, 8).Value = 1.152559531239945
$ws.Cells.Item(25, 9).Value = 1.141925222516107
$ws.Cells.Item(25, 11).Value = 0.4441438060311214
$ws.Cells.Item(25, 12).Value = 0.2281335332629055
